$wb = $excel.ActiveWorkbook

# Delete Sheet4 through Sheet12, leaving Sheet1, Sheet2, Sheet3
$excel.DisplayAlerts = $false
for ($i = 12; $i -ge 4; $i--) {
    $wb.Worksheets.Item("Sheet$i").Delete()
}

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Put a single space in A1 of Sheet2 and Sheet3
$ws2.Range("A1").Value = " "
$ws3.Range("A1").Value = " "

# Select cells
$ws1.Range("A2").Select()
$ws3.Range("A6").Select()

# Activate Sheet3 last (becomes the active tab)
$ws3.Activate()
